$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "5216421-65.2022.8.21.0001"
$ws.Range("C2").Value = "5042179-98.2020.8.21.0001"
$ws.Range("D2").Value = "CIV.36110.01"
$ws.Range("E2").Value = "originario_principal"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "5006597-19.2022.8.21.0049"
$ws.Range("C3").Value = "5002821-79.2020.8.21.0049"
$ws.Range("D3").Value = "CIV.36217.01"
$ws.Range("E3").Value = "originario_principal"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "5002093-18.2019.8.21.0067"
$ws.Range("C4").Value = "9000698-20.2019.8.21.0067"
$ws.Range("D4").Value = "CIV.13144.01"
$ws.Range("E4").Value = "originario_principal"

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "5002094-03.2019.8.21.0067"
$ws.Range("C5").Value = "9001132-09.2019.8.21.0067"
$ws.Range("D5").Value = "CIV.11464.01"
$ws.Range("E5").Value = "originario_principal"

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "5000997-02.2018.8.21.0067"
$ws.Range("C6").Value = "9000692-47.2018.8.21.0067"
$ws.Range("D6").Value = "CIV.10955.01"
$ws.Range("E6").Value = "originario_principal"

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "5002039-18.2020.8.21.0067"
$ws.Range("C7").Value = "9000407-83.2020.8.21.0067"
$ws.Range("D7").Value = "CIV.36418.01"
$ws.Range("E7").Value = "originario_principal"

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "5000998-84.2018.8.21.0067"
$ws.Range("C8").Value = "9000870-93.2018.8.21.0067"
$ws.Range("D8").Value = "CIV.06219.01"
$ws.Range("E8").Value = "originario_principal"

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "5002095-85.2019.8.21.0067"
$ws.Range("C9").Value = "9000724-18.2019.8.21.0067"
$ws.Range("D9").Value = "CIV.12347.01"
$ws.Range("E9").Value = "originario_principal"

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "5002096-70.2019.8.21.0067"
$ws.Range("C10").Value = "9000532-85.2019.8.21.0067"
$ws.Range("D10").Value = "CIV.11868.01"
$ws.Range("E10").Value = "originario_principal"

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "5002097-55.2019.8.21.0067"
$ws.Range("C11").Value = "9000464-38.2019.8.21.0067"
$ws.Range("D11").Value = "CIV.06093.01"
$ws.Range("E11").Value = "originario_principal"

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "5053427-32.2018.8.21.0001"
$ws.Range("C12").Value = "0046017-06.2020.8.21.9000"
$ws.Range("D12").Value = "CIV.35035.02"
$ws.Range("E12").Value = "originario_principal"

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "5035085-12.2014.8.21.0001"
$ws.Range("C13").Value = "0395366-43.2014.8.21.0001"
$ws.Range("D13").Value = "CIV.06910.01"
$ws.Range("E13").Value = "originario_principal"

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "5013439-53.2008.8.21.0001"
$ws.Range("C14").Value = "1253081-20.2008.8.21.0001"
$ws.Range("D14").Value = "CIV.15768.01"
$ws.Range("E14").Value = "originario_principal"

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "5013439-53.2008.8.21.0001"
$ws.Range("C15").Value = "1253081-20.2008.8.21.0001"
$ws.Range("D15").Value = "CIV.15768.01"
$ws.Range("E15").Value = "originario_principal"

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "5013439-53.2008.8.21.0001"
$ws.Range("C16").Value = "1253081-20.2008.8.21.0001"
$ws.Range("D16").Value = "CIV.15768.01"
$ws.Range("E16").Value = "originario_principal"

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "5013439-53.2008.8.21.0001"
$ws.Range("C17").Value = "1253081-20.2008.8.21.0001"
$ws.Range("D17").Value = "CIV.15768.01"
$ws.Range("E17").Value = "originario_principal"

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "5013439-53.2008.8.21.0001"
$ws.Range("C18").Value = "1253081-20.2008.8.21.0001"
$ws.Range("D18").Value = "CIV.15768.01"
$ws.Range("E18").Value = "originario_principal"

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "5014763-10.2010.8.21.0001"
$ws.Range("C19").Value = "3103141-73.2010.8.21.0001"
$ws.Range("D19").Value = "CIV.05499.01"
$ws.Range("E19").Value = "originario_principal"

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "5014763-10.2010.8.21.0001"
$ws.Range("C20").Value = "3103141-73.2010.8.21.0001"
$ws.Range("D20").Value = "CIV.05499.01"
$ws.Range("E20").Value = "originario_principal"

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "5014763-10.2010.8.21.0001"
$ws.Range("C21").Value = "3103141-73.2010.8.21.0001"
$ws.Range("D21").Value = "CIV.05499.01"
$ws.Range("E21").Value = "originario_principal"

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "5029278-74.2015.8.21.0001"
$ws.Range("C22").Value = "0199864-35.2015.8.21.0001"
$ws.Range("D22").Value = "CIV.18395.01"
$ws.Range("E22").Value = "originario_principal"

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "5020879-61.2012.8.21.0001"
$ws.Range("C23").Value = "0026384-42.2012.8.21.3001"
$ws.Range("D23").Value = "CIV.01237.01"
$ws.Range("E23").Value = "originario_principal"

# Apply the existing column-A style (bold, bordered, centered) to the newly added rows
$ws.Range("A6").Copy()
$ws.Range("A7:A23").PasteSpecial(-4122)
$excel.CutCopyMode = 0
